$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix exploration time & agent_step_time formula:
# Avg_Agent_Step_Time (G) and Avg_Experiment_Time (H), along with their
# std-dev counterparts Std_Agent_Step_Time (M) and Std_Experiment_Time (N),
# were being computed with a bug. Update them to the corrected values for
# each data row (rows 2-13).

$ws.Range("G2").Value = 9.873367259999998
$ws.Range("H2").Value = 554.20258152
$ws.Range("M2").Value = 0.7670311963579219
$ws.Range("N2").Value = 75.13766751988602

$ws.Range("G3").Value = 11.45747308
$ws.Range("H3").Value = 1042.4603555
$ws.Range("M3").Value = 1.028636362809473
$ws.Range("N3").Value = 187.1265016891786

$ws.Range("G4").Value = 3.58765748
$ws.Range("H4").Value = 102.26213584
$ws.Range("M4").Value = 0.4988750752095963
$ws.Range("N4").Value = 28.54613897638486

$ws.Range("G5").Value = 4.26867462
$ws.Range("H5").Value = 207.35744618
$ws.Range("M5").Value = 0.5337462031227497
$ws.Range("N5").Value = 51.51781421333257

$ws.Range("G6").Value = 1.08146138
$ws.Range("H6").Value = 16.39288652
$ws.Range("M6").Value = 0.27293318143559
$ws.Range("N6").Value = 7.209824996162254

$ws.Range("G7").Value = 1.44500538
$ws.Range("H7").Value = 37.03972221999999
$ws.Range("M7").Value = 0.272269591049158
$ws.Range("N7").Value = 12.71681812016646

$ws.Range("G8").Value = 0.52715958
$ws.Range("H8").Value = 5.2137149
$ws.Range("M8").Value = 0.1780115570982533
$ws.Range("N8").Value = 2.983240074602246

$ws.Range("G9").Value = 0.72010846
$ws.Range("H9").Value = 12.82987868
$ws.Range("M9").Value = 0.1813003426392532
$ws.Range("N9").Value = 6.14235946071109

$ws.Range("G10").Value = 0.28002424
$ws.Range("H10").Value = 2.09416518
$ws.Range("M10").Value = 0.111553048200163
$ws.Range("N10").Value = 1.442618604615973

$ws.Range("G11").Value = 0.41094848
$ws.Range("H11").Value = 5.934013800000001
$ws.Range("M11").Value = 0.1284182967695023
$ws.Range("N11").Value = 3.862474097052401

$ws.Range("G12").Value = 0.16705888
$ws.Range("H12").Value = 1.013773
$ws.Range("M12").Value = 0.07737831648251943
$ws.Range("N12").Value = 0.7380407141602123

$ws.Range("G13").Value = 0.25024638
$ws.Range("H13").Value = 2.94975924
$ws.Range("M13").Value = 0.08711888434268329
$ws.Range("N13").Value = 2.107001015159887
